# Saldo.xlsx update:
#  - two existing accounts (PATRICIA / 004989791 and PEDRO / 005232019) get
#    their balances updated to much larger figures
#  - five new accounts are appended
#  - the whole "Conta/Nome/Saldo" table is re-sorted descending by Saldo,
#    which is how the sheet was ordered before and after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated data row (row 1 is the header; the table is
# immediately followed by a blank row + a "Filtros aplicados" note row).
$xlDown = -4121
$lastDataRow = $ws.Cells(1, 1).End($xlDown).Row

# --- update the two existing accounts whose balance changed -----------------
$patriciaCell = $ws.Cells.Find("004989791")
$ws.Cells($patriciaCell.Row, 3).Value = 39645.44

$pedroCell = $ws.Cells.Find("005232019")
$ws.Cells($pedroCell.Row, 3).Value = 3000.71

# --- append the brand-new accounts ------------------------------------------
$newRows = @(
    @("'004480134", "JOSE", 66000),
    @("'005529100", "DIMITRI", 6000),
    @("'004547722", "MARCIA", 495),
    @("'001761119", "BLUEMETRIX", 214.29),
    @("'005927101", "SIMONE", 30)
)

$insertFirst = $lastDataRow + 1
$insertLast = $lastDataRow + $newRows.Count
$insertAddress = "$insertFirst" + ":" + "$insertLast"
$ws.Rows($insertAddress).Insert()

$r = $lastDataRow + 1
foreach ($row in $newRows) {
    $ws.Cells($r, 1).Value = $row[0]
    $ws.Cells($r, 2).Value = $row[1]
    $ws.Cells($r, 3).Value = $row[2]
    $r = $r + 1
}

$newLastDataRow = $lastDataRow + $newRows.Count

# --- re-sort the full data range (excludes header row 1) descending by Saldo
$xlDescending = 2
$sortRange = $ws.Range($ws.Cells(2, 1), $ws.Cells($newLastDataRow, 3))
$sortKey = $ws.Range($ws.Cells(2, 3), $ws.Cells($newLastDataRow, 3))
$sortRange.Sort($sortKey, $xlDescending)
